# Updates crypto price/volume data and swaps two rows per the
# "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.12"
$ws.Range("E2").Value = "'0.50%"
$ws.Range("D3").Value = "'27.42"
$ws.Range("E3").Value = "'0.64%"
$ws.Range("D4").Value = "'4.839"
$ws.Range("E4").Value = "'0.59%"
$ws.Range("D5").Value = "'0.06380"
$ws.Range("E5").Value = "'0.13%"
$ws.Range("D6").Value = "'7.035"
$ws.Range("E6").Value = "'1.06%"
$ws.Range("D7").Value = "'1.290"
$ws.Range("E7").Value = "'-0.10%"
$ws.Range("D8").Value = "'0.8940"
$ws.Range("E8").Value = "'1.70%"
$ws.Range("D9").Value = "'0.1520"
$ws.Range("E9").Value = "'-0.33%"
$ws.Range("D10").Value = "'0.05839"
$ws.Range("E10").Value = "'15.27%"
$ws.Range("D11").Value = "'0.07489"
$ws.Range("E11").Value = "'-0.33%"
$ws.Range("D12").Value = "'0.02935"
$ws.Range("E12").Value = "'-1.06%"
$ws.Range("D13").Value = "'0.08991"
$ws.Range("E13").Value = "'-0.32%"
$ws.Range("D14").Value = "'0.001566"
$ws.Range("E14").Value = "'0.15%"
$ws.Range("D15").Value = "'0.0006400"
$ws.Range("E15").Value = "'0.13%"
$ws.Range("D16").Value = "'0.006073"
$ws.Range("E16").Value = "'1.54%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.81%"
$ws.Range("D18").Value = "'3.305"
$ws.Range("E18").Value = "'-0.01%"
$ws.Range("D19").Value = "'2.227"
$ws.Range("E19").Value = "'-1.97%"
$ws.Range("E20").Value = "'-0.78%"
$ws.Range("E21").Value = "'1.14%"
$ws.Range("D22").Value = "'3.903"
$ws.Range("E22").Value = "'-0.54%"
$ws.Range("D23").Value = "'0.04405"
$ws.Range("E23").Value = "'-0.13%"
$ws.Range("D24").Value = "'0.1502"
$ws.Range("E24").Value = "'8.85%"
$ws.Range("D25").Value = "'0.001178"
$ws.Range("E25").Value = "'0.53%"
$ws.Range("E26").Value = "'10.53%"
$ws.Range("D29").Value = "'0.0001652"
$ws.Range("E29").Value = "'-14.66%"
$ws.Range("D40").Value = "'0.04080"
$ws.Range("E40").Value = "'-1.52%"
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1407"
$ws.Range("E41").Value = "'19.15%"
$ws.Range("B42").Value = "'KickToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006601"
$ws.Range("E42").Value = "'-3.12%"
$ws.Range("D43").Value = "'0.002130"
$ws.Range("D44").Value = "'0.01095"
$ws.Range("E44").Value = "'-2.40%"
$ws.Range("D45").Value = "'0.00005527"
$ws.Range("E45").Value = "'6.97%"
$ws.Range("D47").Value = "'0.01848"
$ws.Range("E47").Value = "'-8.54%"
